# Fruta / hortaliza, semanal
# Insert a new weekly price-report row for "Haba" (row 56) in the sheet,
# shifting all subsequent rows down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 56; this shifts rows 56:110 down to 57:111
# and extends the used range / dimension to A1:R111 automatically.
$ws.Rows.Item(56).Insert()

# Populate the newly inserted row 56 with the new weekly record.
$ws.Cells.Item(56, 1).Value = 2
$ws.Cells.Item(56, 2).Value = "Comercializadora del Agro de Limarí"
$ws.Cells.Item(56, 3).Value = "Coquimbo"
$ws.Cells.Item(56, 4).Value = "2023-09-06"
$ws.Cells.Item(56, 5).Value = 4
$ws.Cells.Item(56, 6).Value = 100112026
$ws.Cells.Item(56, 7).Value = "Haba"
$ws.Cells.Item(56, 8).Value = "Sin especificar"
$ws.Cells.Item(56, 9).Value = "Primera"
$ws.Cells.Item(56, 10).Value = 1100
$ws.Cells.Item(56, 11).Value = 9000
$ws.Cells.Item(56, 12).Value = 10000
$ws.Cells.Item(56, 13).Value = 9500
$ws.Cells.Item(56, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(56, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(56, 16).Value = 380
$ws.Cells.Item(56, 17).Value = 25
$ws.Cells.Item(56, 18).Value = "Hortaliza"
